# Moving from 3.1.1 to 3.2.0.
# Insert a "M2Doc version mismatch" warning (4-space run + two orange/highlighted
# runs: "<---" and the message) right after the existing stack-trace run and
# before the existing "    demonstration" run.
#
# Word's Range.InsertAfter/InsertBefore always merges new text into whichever
# run is adjacent at a run boundary (inheriting its formatting), so directly
# inserting+(re)formatting at the target location either bleeds the new
# formatting into the neighbouring run or leaves stray rPr attributes behind
# instead of a clean, attribute-free run. To avoid that, the three new runs
# are authored from scratch in a disposable scratch paragraph appended at the
# very end of the document (a context with no inherited formatting at all),
# copied as FormattedText into the real target location (which cleanly
# creates independent runs without disturbing the neighbouring runs), and the
# scratch paragraph is then deleted.

$d = $word.ActiveDocument

# 1. Append a throw-away paragraph at the very end of the document.
$endRng = $d.Content
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()

# 2. Type the three new runs into the scratch paragraph with the exact
#    formatting required by the target.
$scratch = $d.Content
$scratch.Collapse(0)
$scratch.InsertAfter("    ")
$scratchStart = $scratch.Start

$scratch.Collapse(0)
$scratch.InsertAfter("<---")
$scratch.Font.Color = 42495
$scratch.Font.Size = 16
$scratch.Font.HighlightColorIndex = "#C0C0C0"

$scratch.Collapse(0)
$scratch.InsertAfter("M2Doc version mismatch: template is 3.1.1 and runtime is 3.2.0")
$scratch.Font.Color = 42495
$scratch.Font.Size = 16
$scratch.Font.HighlightColorIndex = "#C0C0C0"
$scratchEnd = $scratch.End

# 3. Grab the formatted (multi-run) content of the scratch paragraph.
$fullScratch = $d.Range($scratchStart, $scratchEnd)
$ft = $fullScratch.FormattedText

# 4. Find the end of the stack trace run (including its trailing newline, so
#    the original run is not split) and paste the captured runs right there.
$target = $d.Content
$target.Find.Execute("Thread.run(Thread.java:748)`n", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target.Collapse(0)
$target.FormattedText = $ft

# 5. Remove the now-unneeded scratch paragraph (still the last paragraph).
$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.Delete()
